$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update / clear cell values (rows 10, 13-24 content realignment) ---
$ws.Range("B10").Value = "Fornecer ao estudante os principais tipos de síntese orgânica e inorgânica de materiais bem como apresentar as principais técnicas analíticas para caracterização de materiais."
$ws.Range("C10").Value = "Fornecer ao estudante os principais tipos de síntese orgânica e inorgânica de materiais bem como apresentar as principais técnicas analíticas para caracterização de materiais."
$ws.Range("A13").Clear()
$ws.Range("B13").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C13").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("A14").Clear()
$ws.Range("B14").Value = "5840897 - Clodoaldo Saron"
$ws.Range("C14").Value = "5840897 - Clodoaldo Saron"
$ws.Range("A15").Value = "Programa resumido:"
$ws.Range("B15").Value = "Introdução à química e sua associação com síntese de novos materiais. A visão moderna do átomo  e Ligações químicas. Estrutura cristalina e técnicas de caracterização cristalográfica. Filmes finos epitaxiais e filmes de uma maneira geral e seu impacto na tecnologica moderna. Crescimento de cristais  Materiais amorfos, síntese e aplicações. Processos e Técnicas de crescimento de cristais de um modo geral. Polímeros condutores e suas aplicações em tecnologica moderna."
$ws.Range("C15").Value = "Introdução à química e sua associação com síntese de novos materiais. A visão moderna do átomo  e Ligações químicas. Estrutura cristalina e técnicas de caracterização cristalográfica. Filmes finos epitaxiais e filmes de uma maneira geral e seu impacto na tecnologica moderna. Crescimento de cristais  Materiais amorfos, síntese e aplicações. Processos e Técnicas de crescimento de cristais de um modo geral. Polímeros condutores e suas aplicações em tecnologica moderna."
$ws.Range("A16").Value = "Short syllabus:"
$ws.Range("B16").Value = "Introduction to the chemistry of materials and its association with the synthesis of new materials. The modern view of the atom and chemical bonds. Crystal structure and crystallographic characterization techniques. Epitaxial thin films and films in general and their impact on modern technology. Amorphous materials, synthesis and applications. Synthesis of materials and chemical transformations. Processes and Techniques of crystal growth in general. Conducting polymers and their applications in modern technology."
$ws.Range("C16").Value = "Introduction to the chemistry of materials and its association with the synthesis of new materials. The modern view of the atom and chemical bonds. Crystal structure and crystallographic characterization techniques. Epitaxial thin films and films in general and their impact on modern technology. Amorphous materials, synthesis and applications. Synthesis of materials and chemical transformations. Processes and Techniques of crystal growth in general. Conducting polymers and their applications in modern technology."
$ws.Range("A17").Value = "Programa:"
$ws.Range("B17").Value = "Química de materiais: definição; papel da química na ciência de materiais; fundamentos.Atomística e a visão moderna do átomo com fundamentos quânticos.Tipos de ligações químicas: forças de van der Waals, potencial de Lennard-Jones, ligação covalente, ligações por coordenação, ligações iônicas e ligações metálicas.Materiais policristalinos e monocristalinos. A ordem cristalográfica e técnicas de caracterização cristalográfica e microscópica. A importância de monocristais em aplicações eletrônicas. Técnicas de crescimento de cristais de alta qualidade tais como: método do fluxo, método Czochralski, método Brigdmann, método do transporte de vapor e método de crescimento de transporte de vapor modificado e isotérmico. Materiais amorfos e sua importância para a tecnologica moderna. Conceitos e técnicas de crescimento de materiais amorfos. Filmes finos epitaxiais, técnicas de crescimento tais como: vapor químico, sputtering, laser ablation e MBE. Filmes finos crescidos por eletrólise para revestimento protetivo, conceitos e aplicações. Síntese de polímeros condutores, conceitos e aplicações como dispositivos eletrônicos."
$ws.Range("C17").Value = "Química de materiais: definição; papel da química na ciência de materiais; fundamentos.Atomística e a visão moderna do átomo com fundamentos quânticos.Tipos de ligações químicas: forças de van der Waals, potencial de Lennard-Jones, ligação covalente, ligações por coordenação, ligações iônicas e ligações metálicas.Materiais policristalinos e monocristalinos. A ordem cristalográfica e técnicas de caracterização cristalográfica e microscópica. A importância de monocristais em aplicações eletrônicas. Técnicas de crescimento de cristais de alta qualidade tais como: método do fluxo, método Czochralski, método Brigdmann, método do transporte de vapor e método de crescimento de transporte de vapor modificado e isotérmico. Materiais amorfos e sua importância para a tecnologica moderna. Conceitos e técnicas de crescimento de materiais amorfos. Filmes finos epitaxiais, técnicas de crescimento tais como: vapor químico, sputtering, laser ablation e MBE. Filmes finos crescidos por eletrólise para revestimento protetivo, conceitos e aplicações. Síntese de polímeros condutores, conceitos e aplicações como dispositivos eletrônicos."
$ws.Range("A18").Value = "Syllabus:"
$ws.Range("B18").Value = "Materials chemistry: definition; role of chemistry in materials science; fundamentals.Atomistics and the modern view of the atom with quantum foundations.Types of chemical bonds: van der Waals forces, Lennard-Jones potential, covalent bonding, coordination bonds, ionic bonds and metallic bonds.Polycrystalline and monocrystalline materials. The crystallographic order and crystallographic and microscopic characterization techniques. The importance of single crystals in electronic applications. High quality crystal growth techniques such as: flow method, Czochralski method, Brigdmann method, vapor transport method and modified isothermal vapor transport growth method. Amorphous materials and their importance for modern technology. Concepts and techniques for growing amorphous materials. Epitaxial thin films, growth techniques such as: chemical vapor, sputtering, laser ablation and MBE. Thin films grown by electrolysis for protective coating, concepts and applications. Synthesis of conductive polymers, concepts and applications as electronic devices."
$ws.Range("C18").Value = "Materials chemistry: definition; role of chemistry in materials science; fundamentals.Atomistics and the modern view of the atom with quantum foundations.Types of chemical bonds: van der Waals forces, Lennard-Jones potential, covalent bonding, coordination bonds, ionic bonds and metallic bonds.Polycrystalline and monocrystalline materials. The crystallographic order and crystallographic and microscopic characterization techniques. The importance of single crystals in electronic applications. High quality crystal growth techniques such as: flow method, Czochralski method, Brigdmann method, vapor transport method and modified isothermal vapor transport growth method. Amorphous materials and their importance for modern technology. Concepts and techniques for growing amorphous materials. Epitaxial thin films, growth techniques such as: chemical vapor, sputtering, laser ablation and MBE. Thin films grown by electrolysis for protective coating, concepts and applications. Synthesis of conductive polymers, concepts and applications as electronic devices."
$ws.Range("A19").Value = "Avaliação:"
$ws.Range("B19").Clear()
$ws.Range("C19").Clear()
$ws.Range("A20").Value = "Método:"
$ws.Range("B20").Value = "Aulas expositivas e práticas ministradas em laboratório."
$ws.Range("C20").Value = "Aulas expositivas e práticas ministradas em laboratório."
$ws.Range("A21").Value = "Critério:"
$ws.Range("B21").Value = "Média simples de duas provas escritas,  Conceito Final = (P1 + P2)/2"
$ws.Range("C21").Value = "Média simples de duas provas escritas,  Conceito Final = (P1 + P2)/2"
$ws.Range("A22").Value = "Norma de recuperação:"
$ws.Range("B22").Value = "Aplicação de duas provas escritas dentro do prazo regimental antes do início do próximo semestre letivo."
$ws.Range("C22").Value = "Aplicação de duas provas escritas dentro do prazo regimental antes do início do próximo semestre letivo."
$ws.Range("A23").Value = "Bibliografia:"
$ws.Range("B23").Value = "ALLCOCK, H. R. Introduction to Materials Chemistry, Wiley, Nova Iorque, 2008.`nFAHLMAN, B. D. Materials Chemistry, Springer, Holanda, 2007.`nZHANG, S.; LI, L.; KUMAR, A. Materials Characterization Techniques, Boca Raton: CRC Press, 2008.`nLENG, Y. Materials Characterization: Introduction to Microscopic and Spectroscopic Methods, Wiley, Cingapura, 2008."
$ws.Range("C23").Value = "ALLCOCK, H. R. Introduction to Materials Chemistry, Wiley, Nova Iorque, 2008.`nFAHLMAN, B. D. Materials Chemistry, Springer, Holanda, 2007.`nZHANG, S.; LI, L.; KUMAR, A. Materials Characterization Techniques, Boca Raton: CRC Press, 2008.`nLENG, Y. Materials Characterization: Introduction to Microscopic and Spectroscopic Methods, Wiley, Cingapura, 2008."
$ws.Range("A24").Value = "Requisitos:"
$ws.Range("B25").Value = "LOM3240 -  Química Inorgânica Fundamental e Aplicada  (Requisito)`n"
$ws.Range("C25").Value = "LOM3240 -  Química Inorgânica Fundamental e Aplicada  (Requisito)`n"

# --- Fix styles for brand-new cells in column B (avoid inheriting column A style) ---
# Column B/C defaults only resolve correctly for pre-existing cells; newly created
# cells in these rows must copy number/font formatting from an existing B/C cell.
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4122) | Out-Null
$ws.Range("C10").Copy() | Out-Null
$ws.Range("C17").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B22").PasteSpecial(-4122) | Out-Null
$ws.Range("C10").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B25").PasteSpecial(-4122) | Out-Null
$ws.Range("C10").Copy() | Out-Null
$ws.Range("C25").PasteSpecial(-4122) | Out-Null

# --- Row heights ---
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(14).AutoFit()
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 120
$ws.Rows.Item(19).AutoFit()
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(23).RowHeight = 120
$ws.Rows.Item(25).RowHeight = 30
